# Insert a new data row at row 421 (pushing the existing rows 421-450 down to
# 422-451) on the single worksheet of the workbook, then populate the new
# row with its values. This mirrors the diff, where every existing row from
# 421 onward shifted down by one and a brand-new row of data was added at
# the top of that block (with the dimension growing from R450 to R451).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 421:450 down to 422:451, leaving row 421 empty.
$ws.Rows.Item(421).Insert()

# Populate the newly inserted row 421 with its new data.
$ws.Range("A421").Value = 9
$ws.Range("B421").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C421").Value = "Metropolitana"
$ws.Range("D421").Value = 44746
$ws.Range("E421").Value = 13
$ws.Range("F421").Value = 100112012
$ws.Range("G421").Value = "Espinaca"
$ws.Range("H421").Value = "Sin especificar"
$ws.Range("I421").Value = "Primera"
$ws.Range("J421").Value = 70
$ws.Range("K421").Value = 9000
$ws.Range("L421").Value = 10000
$ws.Range("M421").Value = 9500
$ws.Range("N421").Value = "`$/cuna 10 kilos"
$ws.Range("O421").Value = "Provincia de Chacabuco"
$ws.Range("P421").Value = 950
$ws.Range("Q421").Value = 10
$ws.Range("R421").Value = "Hortaliza"
